$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence being revised.
$pr = $null
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Our objective is to analyze*") {
        $pr = $p.Range
        break
    }
}

if ($pr -eq $null) {
    throw "Could not locate target paragraph ('Our objective is to analyze ...')."
}

# Common run properties (Arial, bold) shared by every run in this paragraph.
$rpr = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr>'

$runFirst  = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr><w:t>In the fourth quarter of 2023, the US stock market is experiencing significant volatility</w:t></w:r>'
$run01     = '<w:r w:rsidR="00D40C19">' + $rpr + '<w:t xml:space="preserve">. </w:t></w:r>'
$run02     = '<w:r>' + $rpr + '<w:t>O</w:t></w:r>'
$run03     = '<w:r>' + $rpr + '<w:t>bjective is to analyze and understand the current environment across multiple dimensions in order to predict with sound data science principles where the US stock market will be by the end of the first quarter of 2024.</w:t></w:r>'
$run04     = '<w:r>' + $rpr + '<w:t xml:space="preserve"> Success criteria includes whether </w:t></w:r>'
$run05     = '<w:r>' + $rpr + '<w:t xml:space="preserve">three-month </w:t></w:r>'
$run06     = '<w:r>' + $rpr + '<w:t xml:space="preserve">prediction </w:t></w:r>'
$run07     = '<w:r>' + $rpr + '<w:t>interval of</w:t></w:r>'
$run08     = '<w:r>' + $rpr + '<w:t xml:space="preserve"> S&amp;P 500 was </w:t></w:r>'
$run09     = '<w:r>' + $rpr + '<w:t>within</w:t></w:r>'
$run10     = '<w:r>' + $rpr + '<w:t xml:space="preserve"> within a 90% confidence interval</w:t></w:r>'
$run11     = '<w:r>' + $rpr + '<w:t xml:space="preserve"> for at least 90% of periods and final </w:t></w:r>'
$run12     = '<w:r>' + $rpr + '<w:t xml:space="preserve">price </w:t></w:r>'
$run13     = '<w:r>' + $rpr + '<w:t xml:space="preserve">APE &lt; </w:t></w:r>'
$run14     = '<w:r>' + $rpr + '<w:t>1</w:t></w:r>'
$run15     = '<w:r>' + $rpr + '<w:t>0%.</w:t></w:r>'

$pPr = '<w:pPr><w:tabs><w:tab w:val="left" w:pos="7020"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr>'

$paraOpen = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2252F81B" w14:textId="77777777" w:rsidR="008308C4" w:rsidRDefault="0081339E">' + $pPr

$xmlFrag = $paraOpen + $runFirst + $run01 + $run02 + $run03 + $run04 + $run05 + $run06 + $run07 + $run08 + $run09 + $run10 + $run11 + $run12 + $run13 + $run14 + $run15 + '</w:p>'

# InsertXML, called on a Range spanning the *entire* paragraph, replaces that
# paragraph's contents with the supplied OOXML (preserving paragraph identity).
$null = $pr.InsertXML($xmlFrag)
